$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "category" value for the testProgram row (row 2, column D) to
# the existing "Test Program" string, replicating the new data cell added
# in the diff.
$ws.Range("D2").Value = "Test Program"

# Move/confirm the active selection to D2, matching the saved selection
# state in the workbook.
$ws.Range("D2").Select()
